$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose displayed value looks numeric/percentage must be forced back to
# plain text (matching the source data, which stores these as literal strings)
# via a leading apostrophe, then the quote-prefix style that introduces is reset
# back to "Normal" so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "292.87"
Set-TextValue "E2" "-5.14%"
Set-TextValue "D3" "40.25"
Set-TextValue "E3" "-1.95%"
Set-TextValue "D4" "5.033"
Set-TextValue "E4" "-3.36%"
Set-TextValue "D5" "0.07410"
Set-TextValue "E5" "-3.56%"
Set-TextValue "D6" "1.574"
Set-TextValue "E6" "-3.22%"
Set-TextValue "D7" "0.9313"
Set-TextValue "E7" "1.80%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D8" "2.420"
Set-TextValue "E8" "-0.98%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1190"
Set-TextValue "E9" "-1.88%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1746"
Set-TextValue "E10" "-4.07%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.08767"
Set-TextValue "E11" "-3.90%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.04186"
Set-TextValue "E12" "-0.70%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.1050"
Set-TextValue "E13" "-0.14%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001266"
Set-TextValue "E14" "0.22%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005869"
Set-TextValue "E15" "-0.03%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.361"
Set-TextValue "E16" "0.58%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.342"
Set-TextValue "E17" "0.97%"
Set-TextValue "D18" "0.3350"
Set-TextValue "E18" "0.43%"
Set-TextValue "D19" "7.708"
Set-TextValue "D20" "0.1365"
Set-TextValue "E20" "-1.27%"
Set-TextValue "D21" "0.2824"
Set-TextValue "E21" "4.16%"
Set-TextValue "D22" "0.03873"
Set-TextValue "E22" "-3.50%"
Set-TextValue "D23" "0.001297"
Set-TextValue "E23" "2.88%"
Set-TextValue "D24" "0.003484"
Set-TextValue "E24" "-18.44%"
Set-TextValue "D25" "0.0001309"
Set-TextValue "E25" "0.67%"
Set-TextValue "D26" "0.0003749"
Set-TextValue "D38" "0.02297"
Set-TextValue "E38" "-8.16%"
Set-TextValue "D39" "0.04995"
Set-TextValue "E39" "-6.04%"
Set-TextValue "D40" "0.007698"
Set-TextValue "E40" "-1.69%"
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D41" "0.004169"
Set-TextValue "E41" "125.05%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1276"
Set-TextValue "E42" "-2.71%"
Set-TextValue "D43" "0.007503"
Set-TextValue "E43" "15.35%"
Set-TextValue "D44" "0.007157"
Set-TextValue "E44" "-12.69%"
Set-TextValue "D45" "0.3184"
Set-TextValue "E45" "-4.68%"
Set-TextValue "D46" "0.00006742"
Set-TextValue "E46" "0.44%"
Set-TextValue "D47" "0.00000000755"
Set-TextValue "E47" "0.70%"
Set-TextValue "D48" "0.2518"
Set-TextValue "E48" "-16.82%"
Set-TextValue "D49" "0.004230"
Set-TextValue "E49" "36.43%"
Set-TextValue "D50" "0.00002115"
Set-TextValue "E50" "0.70%"
Set-TextValue "D51" "0.0002014"
Set-TextValue "E51" "0.70%"
